# Scheduled-runner price/profit refresh across the Leve tracking sheets.
# Updates currentAveragePrice / LevePrice / LeveProfit columns (H-N) for the
# rows whose market data changed since the last pull. A few rows gain or
# lose trailing HQ columns (M/N) depending on whether an HQ price exists.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1211.3529
$ws.Range("I15").Value = 1211.3529
$ws.Range("K15").Value = 3634.0587
$ws.Range("M15").Value = -3465.0587

$ws.Range("H33").Value = 231.5625
$ws.Range("I33").Value = 193.66667
$ws.Range("K33").Value = 193.66667
$ws.Range("M33").Value = 35.33332999999999

$ws.Range("H43").Value = 7350
$ws.Range("J43").Value = 7350
$ws.Range("L43").Value = 7350
$ws.Range("N43").Value = -7488

$ws.Range("H62").Value = 6315.6665
$ws.Range("I62").Value = 5968
$ws.Range("K62").Value = 5968
$ws.Range("M62").Value = -5344

$ws.Range("H65").Value = 6315.6665
$ws.Range("I65").Value = 5968
$ws.Range("K65").Value = 29840
$ws.Range("M65").Value = -26720

$ws.Range("H101").Value = 0
$ws.Range("I101").Value = 0
$ws.Range("K101").Value = 0
$ws.Range("M101").Value = $null

$ws.Range("H107").Value = 1322.5
$ws.Range("I107").Value = 1322.5
$ws.Range("K107").Value = 1322.5
$ws.Range("M107").Value = 597.5

$ws.Range("H127").Value = 3586.3333
$ws.Range("I127").Value = 379.5
$ws.Range("J127").Value = 10000
$ws.Range("K127").Value = 1138.5
$ws.Range("L127").Value = 30000
$ws.Range("M127").Value = 3821.5
$ws.Range("N127").Value = -39920

$ws.Range("H131").Value = 8499
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 8499
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 25497
$ws.Range("M131").Value = $null
$ws.Range("N131").Value = -35577

$ws.Range("H132").Value = 1714.8572
$ws.Range("I132").Value = 2098.6
$ws.Range("K132").Value = 6295.799999999999
$ws.Range("M132").Value = -3765.799999999999

$ws.Range("H137").Value = 4694.1177
$ws.Range("I137").Value = 2449.5
$ws.Range("K137").Value = 7348.5
$ws.Range("M137").Value = -4798.5

$ws.Range("H138").Value = 4970.8945
$ws.Range("J138").Value = 5291.706
$ws.Range("L138").Value = 15875.118
$ws.Range("N138").Value = -26155.118

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 567.63635
$ws.Range("I2").Value = 488.22223
$ws.Range("J2").Value = 925
$ws.Range("K2").Value = 488.22223
$ws.Range("L2").Value = 925
$ws.Range("M2").Value = -375.22223
$ws.Range("N2").Value = -1151

$ws.Range("H19").Value = 963.3333
$ws.Range("I19").Value = 963.3333
$ws.Range("K19").Value = 963.3333
$ws.Range("M19").Value = -734.3333

$ws.Range("H32").Value = 1875.4259
$ws.Range("I32").Value = 1437.9807
$ws.Range("K32").Value = 1437.9807
$ws.Range("M32").Value = -1150.9807

$ws.Range("H45").Value = 1178.7
$ws.Range("I45").Value = 1166.875
$ws.Range("J45").Value = 1226
$ws.Range("K45").Value = 1166.875
$ws.Range("L45").Value = 1226
$ws.Range("M45").Value = -789.875
$ws.Range("N45").Value = -1980

$ws.Range("H116").Value = 567.63635
$ws.Range("I116").Value = 488.22223
$ws.Range("J116").Value = 925
$ws.Range("K116").Value = 488.22223
$ws.Range("L116").Value = 925
$ws.Range("M116").Value = 1805.77777
$ws.Range("N116").Value = -5513

$ws.Range("H122").Value = 3183.7646
$ws.Range("I122").Value = 3018.5386
$ws.Range("K122").Value = 9055.6158
$ws.Range("M122").Value = -6605.6158

$ws.Range("H139").Value = 125000
$ws.Range("J139").Value = 125000
$ws.Range("L139").Value = 125000
$ws.Range("N139").Value = -135280

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 567.63635
$ws.Range("I3").Value = 488.22223
$ws.Range("J3").Value = 925
$ws.Range("K3").Value = 488.22223
$ws.Range("L3").Value = 925
$ws.Range("M3").Value = -374.22223
$ws.Range("N3").Value = -1153

$ws.Range("H94").Value = 8888
$ws.Range("I94").Value = 0
$ws.Range("K94").Value = 0
$ws.Range("M94").Value = $null

$ws.Range("H105").Value = 3499
$ws.Range("I105").Value = 3498
$ws.Range("K105").Value = 3498
$ws.Range("M105").Value = -1751

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5655.9287
$ws.Range("I31").Value = 2699.2727
$ws.Range("J31").Value = 16497
$ws.Range("K31").Value = 2699.2727
$ws.Range("L31").Value = 16497
$ws.Range("M31").Value = -2404.2727
$ws.Range("N31").Value = -17087

$ws.Range("H34").Value = 5655.9287
$ws.Range("I34").Value = 2699.2727
$ws.Range("J34").Value = 16497
$ws.Range("K34").Value = 2699.2727
$ws.Range("L34").Value = 16497
$ws.Range("M34").Value = -2497.2727
$ws.Range("N34").Value = -16901

$ws.Range("H58").Value = 1750

$ws.Range("H105").Value = 1589
$ws.Range("I105").Value = 1589
$ws.Range("K105").Value = 1589
$ws.Range("M105").Value = 158

$ws.Range("H132").Value = 4460.25
$ws.Range("I132").Value = 2998.3333
$ws.Range("K132").Value = 8994.999899999999
$ws.Range("M132").Value = -6464.999899999999

$ws.Range("H136").Value = 1750

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 3784.7144
$ws.Range("J132").Value = 4082.1667
$ws.Range("L132").Value = 36739.5003
$ws.Range("N132").Value = -41799.5003

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2075.5334
$ws.Range("I102").Value = 1386.1666
$ws.Range("K102").Value = 1386.1666
$ws.Range("M102").Value = 235.8334

$ws.Range("H122").Value = 1512.25
$ws.Range("I122").Value = 1574.75
$ws.Range("K122").Value = 4724.25
$ws.Range("M122").Value = -2274.25

$ws.Range("H126").Value = 1999.6666
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 1999.6666
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 5998.9998
$ws.Range("M126").Value = $null
$ws.Range("N126").Value = -10938.9998

$ws.Range("H132").Value = 3056.5908
$ws.Range("I132").Value = 2802.4736
$ws.Range("K132").Value = 8407.4208
$ws.Range("M132").Value = -5877.4208

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4500
$ws.Range("I7").Value = 4500
$ws.Range("K7").Value = 4500
$ws.Range("M7").Value = -4388

$ws.Range("H93").Value = 2334.3333
$ws.Range("I93").Value = 2499.5
$ws.Range("J93").Value = 2004
$ws.Range("K93").Value = 2499.5
$ws.Range("L93").Value = 2004
$ws.Range("M93").Value = -1251.5
$ws.Range("N93").Value = -4500

$ws.Range("H122").Value = 2221.5
$ws.Range("J122").Value = 2221.5
$ws.Range("L122").Value = 6664.5
$ws.Range("N122").Value = -11564.5

$ws.Range("H126").Value = 4500
$ws.Range("I126").Value = 4500
$ws.Range("K126").Value = 13500
$ws.Range("M126").Value = -11030

$ws.Range("H136").Value = 10744.071
$ws.Range("I136").Value = 11368.583
$ws.Range("K136").Value = 34105.749
$ws.Range("M136").Value = -31555.749

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("N46").Value = $null

$ws.Range("H113").Value = 437.1
$ws.Range("I113").Value = 402.33334
$ws.Range("J113").Value = 750
$ws.Range("K113").Value = 1207.00002
$ws.Range("L113").Value = 2250
$ws.Range("M113").Value = 962.9999800000001
$ws.Range("N113").Value = -6590

$ws.Range("H122").Value = 2500
$ws.Range("I122").Value = 2500
$ws.Range("K122").Value = 7500
$ws.Range("M122").Value = -5050

$ws.Range("H132").Value = 2215.1738
$ws.Range("I132").Value = 1636
$ws.Range("K132").Value = 4908
$ws.Range("M132").Value = -2378

$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").Value = $null

$ws.Range("H136").Value = 7947.476
$ws.Range("I136").Value = 7772
$ws.Range("K136").Value = 23316
$ws.Range("M136").Value = -20766
